$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (shifts old B->C, old C->D)
$ws.Columns.Item(2).Insert()

# New column B width should match column A's width (~75.8 characters)
$ws.Columns.Item(2).ColumnWidth = 75

# Header for new column B
$ws.Range("B1").Value = "StatQuery"

# New StatQuery text in B2, matching the wrap-text style used by A2
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE diag.disease_term IN['Stage 3']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").WrapText = $true

# Restore selection to A2 (mirrors the saved sheet view)
$ws.Range("A2").Select()
